$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.723.16'
$ws.Range('E2').Value = '  -6.52%  '
$ws.Range('D3').Value = '2.540.05'
$ws.Range('E3').Value = '  -4.87%  '
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '298.86'
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('D6').Value = '93.53'
$ws.Range('E6').Value = '  -5.09%  '
$ws.Range('D7').Value = '0.573'
$ws.Range('E7').Value = '  -4.37%  '
$ws.Range('E8').Value = '  -0.10%  '
$ws.Range('D9').Value = '0.546'
$ws.Range('E9').Value = '  -6.33%  '
$ws.Range('D10').Value = '35.87'
$ws.Range('E10').Value = '  -6.84%  '
$ws.Range('D11').Value = '0.0804'
$ws.Range('E11').Value = '  -5.50%  '
$ws.Range('D12').Value = '7.69'
$ws.Range('E12').Value = '  -5.69%  '
$ws.Range('E13').Value = '  +6.12%  '
$ws.Range('D14').Value = '2.925.98'
$ws.Range('E14').Value = '  -5.17%  '
$ws.Range('D15').Value = '2.546.46'
$ws.Range('E15').Value = '  -4.36%  '
$ws.Range('D16').Value = '0.876'
$ws.Range('E16').Value = '  -6.61%  '
$ws.Range('D17').Value = '14.18'
$ws.Range('E17').Value = '  -6.70%  '
$ws.Range('D18').Value = '42.710.16'
$ws.Range('E18').Value = '  -6.72%  '
$ws.Range('D19').Value = '0.0₃0979'
$ws.Range('E19').Value = '  -4.54%  '
$ws.Range('E20').Value = '  -2.21%  '
$ws.Range('E21').Value = '  -4.53%  '
$ws.Range('D22').Value = '71.68'
$ws.Range('E22').Value = '  -5.04%  '
$ws.Range('D23').Value = '254.40'
$ws.Range('E23').Value = '  -10.33%  '
$ws.Range('E24').Value = '  -5.39%  '
$ws.Range('E25').Value = '  -5.44%  '
$ws.Range('D26').Value = '29.01'
$ws.Range('E26').Value = '  -6.21%  '
$ws.Range('E27').Value = '  +0.24%  '
$ws.Range('D28').Value = '10.15'
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('D29').Value = '2.11'
$ws.Range('E29').Value = '  -5.04%  '
$ws.Range('D30').Value = '36.52'
$ws.Range('E30').Value = '  -5.61%  '
$ws.Range('D31').Value = '6.06'
$ws.Range('E31').Value = '  -2.61%  '
$ws.Range('D32').Value = '152.51'
$ws.Range('D33').Value = '2.75'
$ws.Range('E33').Value = '  -2.57%  '
$ws.Range('E34').Value = '  -8.95%  '
$ws.Range('E35').Value = '  -10.34%  '
$ws.Range('D36').Value = '0.0792'
$ws.Range('E36').Value = '  -6.31%  '
$ws.Range('E37').Value = '  -6.86%  '
$ws.Range('B38').Value = 'Celestia'
$ws.Range('C38').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D38').Value = '17.23'
$ws.Range('E38').Value = '  +5.25%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').Value = '0.119'
$ws.Range('E39').Value = '  -4.35%  '
$ws.Range('D40').Value = '23.17'
$ws.Range('E40').Value = '  -10.23%  '
$ws.Range('D41').Value = '0.0309'
$ws.Range('E41').Value = '  -5.89%  '
$ws.Range('D42').Value = '3.40'
$ws.Range('E42').Value = '  -6.70%  '
$ws.Range('E43').Value = '  -4.66%  '
$ws.Range('D44').Value = '2.081.82'
$ws.Range('E44').Value = '  -1.53%  '
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.15%  '
$ws.Range('D46').Value = '1.77'
$ws.Range('E46').Value = '  +11.47%  '
$ws.Range('D47').Value = '9.06'
$ws.Range('E47').Value = '  -2.67%  '
$ws.Range('D48').Value = '84.16'
$ws.Range('E48').Value = '  -10.48%  '
$ws.Range('D49').Value = '2.783.92'
$ws.Range('E49').Value = '  -5.05%  '
$ws.Range('D50').Value = '104.59'
$ws.Range('E50').Value = '  -6.75%  '
